$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 611
$ws.Range("F3").Value = 281
$ws.Range("G3").Value = 65
$ws.Range("F5").Value = 760
$ws.Range("F6").Value = 412
$ws.Range("F8").Value = 188
$ws.Range("F10").Value = 253
$ws.Range("F11").Value = 6980
$ws.Range("F12").Value = 71
$ws.Range("F13").Value = 58
$ws.Range("F15").Value = 556
$ws.Range("F16").Value = 377
$ws.Range("F20").Value = 726
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 190
$ws.Range("F23").Value = 106
$ws.Range("F24").Value = 337
$ws.Range("F25").Value = 1046
$ws.Range("F27").Value = 22
$ws.Range("F28").Value = 1957
$ws.Range("F29").Value = 543
$ws.Range("F31").Value = 539

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 304
$ws.Range("F10").Value = 136

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 317

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 317
$ws.Range("F3").Value = 611
$ws.Range("F4").Value = 281
$ws.Range("G4").Value = 65
$ws.Range("F6").Value = 760
$ws.Range("F8").Value = 412
$ws.Range("F10").Value = 188
$ws.Range("F12").Value = 253
$ws.Range("F13").Value = 6980
$ws.Range("F14").Value = 71
$ws.Range("F15").Value = 58
$ws.Range("F18").Value = 556
$ws.Range("F19").Value = 377
$ws.Range("F25").Value = 304
$ws.Range("F27").Value = 726
$ws.Range("F29").Value = 1
$ws.Range("F31").Value = 136
$ws.Range("F32").Value = 190
$ws.Range("F33").Value = 106
$ws.Range("F34").Value = 337
$ws.Range("F35").Value = 1046
$ws.Range("F37").Value = 22
$ws.Range("F38").Value = 1957
$ws.Range("F39").Value = 543
$ws.Range("F41").Value = 539
